# Auto commit at 2025-11-06 12:17:41.39
# Append two new daily rows (date serial 45966 == 2025-11-05) for both
# charging stations to the bottom of the existing data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 72: 四方坪站 (same shared string used by every prior odd-numbered
# data row, e.g. A2/A4/...).
$ws.Cells.Item(72, 1).Value = 45966
$ws.Cells.Item(72, 2).Value = "四方坪站"
$ws.Cells.Item(72, 3).Value = 7953.77
$ws.Cells.Item(72, 4).Value = 6903.5
$ws.Cells.Item(72, 5).Value = 2678.68
$ws.Cells.Item(72, 6).Value = 360

# Row 73: 高岭站
$ws.Cells.Item(73, 1).Value = 45966
$ws.Cells.Item(73, 2).Value = "高岭站"
$ws.Cells.Item(73, 3).Value = 5558.54
$ws.Cells.Item(73, 4).Value = 4987.96
$ws.Cells.Item(73, 5).Value = 1415.89
$ws.Cells.Item(73, 6).Value = 228

# Move the active selection to mirror Excel's usual "select just past the
# new last data row" behaviour after appending rows.
$ws.Range("H71").Select()
